{"js": "const replacements = [\n  [\"2024-04-06 Saturday\", \"2024-04-07 Sunday\"],\n  [\"335\u00f73=111, 2\", \"650\u00f77=92, 6\"],\n  [\"569\u00f76=94, 5\", \"623\u00f77=89, 0\"],\n  [\"673\u00f73=224, 1\", \"143\u00f73=47, 2\"],\n  [\"853\u00f75=170, 3\", \"290\u00f78=36, 2\"],\n  [\"277\u00f79=30, 7\", \"674\u00f76=112, 2\"],\n  [\"953\u00f77=136, 1\", \"320\u00f76=53, 2\"],\n  [\"326\u00f72=163, 0\", \"199\u00f74=49, 3\"],\n  [\"563\u00f76=93, 5\", \"388\u00f74=97, 0\"],\n  [\"122\u00f73=40, 2\", \"301\u00f78=37, 5\"],\n  [\"441\u00f76=73, 3\", \"872\u00f73=290, 2\"],\n  [\"626\u00f78=78, 2\", \"436\u00f77=62, 2\"],\n  [\"918\u00f74=229, 2\", \"860\u00f73=286, 2\"],\n  [\"305\u00f78=38, 1\", \"702\u00f72=351, 0\"],\n  [\"971\u00f75=194, 1\", \"121\u00f78=15, 1\"],\n  [\"452\u00f72=226, 0\", \"584\u00f78=73, 0\"],\n  [\"775\u00f75=155, 0\", \"662\u00f74=165, 2\"],\n  [\"116\u00f73=38, 2\", \"604\u00f72=302, 0\"],\n  [\"308\u00f78=38, 4\", \"434\u00f77=62, 0\"],\n  [\"314\u00f77=44, 6\", \"630\u00f75=126, 0\"],\n  [\"896\u00f75=179, 1\", \"113\u00f77=16, 1\"],\n  [\"173\u00f76=28, 5\", \"288\u00f76=48, 0\"],\n  [\"119\u00f73=39, 2\", \"602\u00f78=75, 2\"],\n  [\"379\u00f74=94, 3\", \"632\u00f73=210, 2\"],\n  [\"460\u00f74=115, 0\", \"397\u00f72=198, 1\"],\n  [\"421\u00f75=84, 1\", \"620\u00f78=77, 4\"],\n];\n\nconst body = context.document.body;\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(after, \"Replace\");\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n$replacements = @(\n    @(\"2024-04-06 Saturday\", \"2024-04-07 Sunday\"),\n    @(\"335\u00f73=111, 2\", \"650\u00f77=92, 6\"),\n    @(\"569\u00f76=94, 5\", \"623\u00f77=89, 0\"),\n    @(\"673\u00f73=224, 1\", \"143\u00f73=47, 2\"),\n    @(\"853\u00f75=170, 3\", \"290\u00f78=36, 2\"),\n    @(\"277\u00f79=30, 7\", \"674\u00f76=112, 2\"),\n    @(\"953\u00f77=136, 1\", \"320\u00f76=53, 2\"),\n    @(\"326\u00f72=163, 0\", \"199\u00f74=49, 3\"),\n    @(\"563\u00f76=93, 5\", \"388\u00f74=97, 0\"),\n    @(\"122\u00f73=40, 2\", \"301\u00f78=37, 5\"),\n    @(\"441\u00f76=73, 3\", \"872\u00f73=290, 2\"),\n    @(\"626\u00f78=78, 2\", \"436\u00f77=62, 2\"),\n    @(\"918\u00f74=229, 2\", \"860\u00f73=286, 2\"),\n    @(\"305\u00f78=38, 1\", \"702\u00f72=351, 0\"),\n    @(\"971\u00f75=194, 1\", \"121\u00f78=15, 1\"),\n    @(\"452\u00f72=226, 0\", \"584\u00f78=73, 0\"),\n    @(\"775\u00f75=155, 0\", \"662\u00f74=165, 2\"),\n    @(\"116\u00f73=38, 2\", \"604\u00f72=302, 0\"),\n    @(\"308\u00f78=38, 4\", \"434\u00f77=62, 0\"),\n    @(\"314\u00f77=44, 6\", \"630\u00f75=126, 0\"),\n    @(\"896\u00f75=179, 1\", \"113\u00f77=16, 1\"),\n    @(\"173\u00f76=28, 5\", \"288\u00f76=48, 0\"),\n    @(\"119\u00f73=39, 2\", \"602\u00f78=75, 2\"),\n    @(\"379\u00f74=94, 3\", \"632\u00f73=210, 2\"),\n    @(\"460\u00f74=115, 0\", \"397\u00f72=198, 1\"),\n    @(\"421\u00f75=84, 1\", \"620\u00f78=77, 4\"),\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $r = $d.Content\n    $r.Find.ClearFormatting()\n    $r.Find.Replacement.ClearFormatting()\n    [void]$r.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}"}
